$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price / Volume(1h) columns (D, E) store their values as literal text
# (e.g. "309.57", "-4.04%") rather than numbers, matching the scraped-data
# source file. A plain `$cell.Value = "310.54"` would let Excel's
# auto-detection reinterpret the numeric- or percent-looking text as a
# real number (and, for "%", stamp a percentage number format onto the
# cell), which would not match the original text-cell representation.
# Prefixing with an apostrophe forces Excel to keep the entry as literal
# text; resetting .Style back to "Normal" afterwards clears the implicit
# "quote prefix" style Excel applies so the cell's formatting stays
# identical to its untouched neighbors.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "310.54"
Set-TextValue $ws.Range("E2") "-3.53%"
Set-TextValue $ws.Range("D3") "49.62"
Set-TextValue $ws.Range("E3") "2.84%"
Set-TextValue $ws.Range("D4") "5.145"
Set-TextValue $ws.Range("E4") "-2.63%"
Set-TextValue $ws.Range("E5") "-4.26%"
Set-TextValue $ws.Range("D6") "4.537"
Set-TextValue $ws.Range("E6") "-0.76%"
Set-TextValue $ws.Range("D7") "1.374"
Set-TextValue $ws.Range("E7") "14.31%"
Set-TextValue $ws.Range("D8") "1.570"
Set-TextValue $ws.Range("E8") "-4.28%"
Set-TextValue $ws.Range("D9") "0.1218"
Set-TextValue $ws.Range("E9") "-6.31%"
Set-TextValue $ws.Range("D10") "0.2003"
Set-TextValue $ws.Range("E10") "2.89%"
Set-TextValue $ws.Range("D11") "0.04739"
Set-TextValue $ws.Range("E11") "2.12%"
Set-TextValue $ws.Range("D12") "0.09291"
Set-TextValue $ws.Range("E12") "-1.93%"
Set-TextValue $ws.Range("D13") "0.1045"
Set-TextValue $ws.Range("E13") "-0.42%"
Set-TextValue $ws.Range("D14") "0.001262"
Set-TextValue $ws.Range("E14") "-5.69%"
Set-TextValue $ws.Range("D15") "0.005823"
Set-TextValue $ws.Range("E15") "0.37%"
Set-TextValue $ws.Range("E16") "2,021.33%"
Set-TextValue $ws.Range("E17") "-0.20%"
Set-TextValue $ws.Range("D18") "2.438"
Set-TextValue $ws.Range("E18") "0.38%"
Set-TextValue $ws.Range("D20") "7.984"
Set-TextValue $ws.Range("E20") "-1.24%"
Set-TextValue $ws.Range("D21") "0.1364"
Set-TextValue $ws.Range("E21") "-2.60%"
Set-TextValue $ws.Range("E23") "0.91%"
Set-TextValue $ws.Range("D24") "0.001274"
Set-TextValue $ws.Range("E24") "-2.46%"
Set-TextValue $ws.Range("D25") "0.003929"
Set-TextValue $ws.Range("E25") "-7.58%"
Set-TextValue $ws.Range("D26") "0.0001349"
Set-TextValue $ws.Range("E26") "-0.10%"
Set-TextValue $ws.Range("D38") "0.02596"
Set-TextValue $ws.Range("E38") "-4.94%"
Set-TextValue $ws.Range("D39") "0.06168"
Set-TextValue $ws.Range("E39") "5.29%"
Set-TextValue $ws.Range("D40") "0.01112"
Set-TextValue $ws.Range("E40") "76.51%"
Set-TextValue $ws.Range("D41") "0.007932"
Set-TextValue $ws.Range("E41") "3.11%"
Set-TextValue $ws.Range("E42") "-1.45%"
Set-TextValue $ws.Range("D43") "0.008383"
Set-TextValue $ws.Range("E43") "8.64%"
Set-TextValue $ws.Range("D44") "0.008333"
Set-TextValue $ws.Range("E44") "2.85%"
Set-TextValue $ws.Range("D45") "0.3118"
Set-TextValue $ws.Range("E45") "-2.38%"
Set-TextValue $ws.Range("D46") "0.00007640"
Set-TextValue $ws.Range("E46") "8.88%"
Set-TextValue $ws.Range("E47") "-0.09%"
Set-TextValue $ws.Range("D48") "0.05325"
Set-TextValue $ws.Range("E48") "-0.97%"
Set-TextValue $ws.Range("D49") "0.002621"
Set-TextValue $ws.Range("E49") "-34.48%"
Set-TextValue $ws.Range("D50") "0.00002098"
Set-TextValue $ws.Range("E50") "-0.09%"
Set-TextValue $ws.Range("D51") "0.0001998"
Set-TextValue $ws.Range("E51") "-0.09%"
